# Update the "想去人数" (interest count) figures on both the "展览"
# and "全部类型" sheets to reflect the latest scraped output.

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 2263
    3 = 1705
    5 = 1086
    6 = 806
    8 = 5825
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
